$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 260, shifting existing rows 260-345 down to 261-346
$ws.Rows("260:260").Insert()

# Populate the newly inserted row 260 with the new data record
$ws.Cells.Item(260, 1).Value = 6
$ws.Cells.Item(260, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(260, 3).Value = "Metropolitana"
$ws.Cells.Item(260, 4).Value = 45120
$ws.Cells.Item(260, 4).NumberFormat = $ws.Cells.Item(259, 4).NumberFormat
$ws.Cells.Item(260, 5).Value = 13
$ws.Cells.Item(260, 6).Value = 100112029
$ws.Cells.Item(260, 7).Value = "Orégano"
$ws.Cells.Item(260, 8).Value = "Sin especificar"
$ws.Cells.Item(260, 9).Value = "Primera"
$ws.Cells.Item(260, 10).Value = 38
$ws.Cells.Item(260, 11).Value = 20000
$ws.Cells.Item(260, 12).Value = 20000
$ws.Cells.Item(260, 13).Value = 20000
$ws.Cells.Item(260, 14).Value = "$/docena de atados"
$ws.Cells.Item(260, 15).Value = "Región Metropolitana"
$ws.Cells.Item(260, 16).Value = 6667
$ws.Cells.Item(260, 17).Value = 3
$ws.Cells.Item(260, 18).Value = "Hortaliza"
